$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (A2/B2)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 183

# New row 3 (A3/B3) - was previously A3=1/B3=175, now becomes A3=0/B3=115
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 115

# Copy style from A2/A3 (already bold/centered/bordered) down to new rows A4, A5
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New row 4 (A4/B4)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 88

# New row 5 (A5/B5)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 60
